$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 15873747
$ws.Cells.Item(41, 9).Value = 447.6
$ws.Cells.Item(41, 10).Value = 30304018
$ws.Cells.Item(41, 11).Value = 447.6
$ws.Cells.Item(41, 12).Value = 30304018
$ws.Cells.Item(41, 13).Value = -7.600000000000023
$ws.Cells.Item(41, 14).Value = -30304898
$ws.Cells.Item(76, 8).Value = 11117511
$ws.Cells.Item(76, 10).Value = 8000
$ws.Cells.Item(76, 12).Value = 8000
$ws.Cells.Item(76, 14).Value = -8630
$ws.Cells.Item(79, 8).Value = 11117511
$ws.Cells.Item(79, 10).Value = 8000
$ws.Cells.Item(79, 12).Value = 8000
$ws.Cells.Item(79, 14).Value = -10184
$ws.Cells.Item(92, 8).Value = 938.3939
$ws.Cells.Item(92, 9).Value = 223.2069
$ws.Cells.Item(92, 11).Value = 223.2069
$ws.Cells.Item(92, 13).Value = 1024.7931
$ws.Cells.Item(100, 8).Value = 1917.1666
$ws.Cells.Item(100, 9).Value = 2100.6
$ws.Cells.Item(100, 11).Value = 2100.6
$ws.Cells.Item(100, 13).Value = -1559.6
$ws.Cells.Item(112, 8).Value = 4033.7
$ws.Cells.Item(112, 10).Value = 4111.5127
$ws.Cells.Item(112, 12).Value = 12334.5381
$ws.Cells.Item(112, 14).Value = -14550.5381

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 6180.2
$ws.Cells.Item(26, 9).Value = 3450.5
$ws.Cells.Item(26, 10).Value = 8000
$ws.Cells.Item(26, 11).Value = 3450.5
$ws.Cells.Item(26, 12).Value = 8000
$ws.Cells.Item(26, 13).Value = -3120.5
$ws.Cells.Item(26, 14).Value = -8660
$ws.Cells.Item(32, 8).Value = 8737.35
$ws.Cells.Item(32, 9).Value = 5094.382
$ws.Cells.Item(32, 10).Value = 16751.88
$ws.Cells.Item(32, 11).Value = 5094.382
$ws.Cells.Item(32, 12).Value = 16751.88
$ws.Cells.Item(32, 13).Value = -4807.382
$ws.Cells.Item(32, 14).Value = -17325.88
$ws.Cells.Item(45, 8).Value = 11993492
$ws.Cells.Item(45, 9).Value = 23977822
$ws.Cells.Item(45, 10).Value = 9162.666999999999
$ws.Cells.Item(45, 11).Value = 23977822
$ws.Cells.Item(45, 12).Value = 9162.666999999999
$ws.Cells.Item(45, 13).Value = -23977445
$ws.Cells.Item(45, 14).Value = -9916.666999999999
$ws.Cells.Item(74, 8).Value = 97206.63
$ws.Cells.Item(74, 9).Value = 102166.445
$ws.Cells.Item(74, 10).Value = 92742.8
$ws.Cells.Item(74, 11).Value = 102166.445
$ws.Cells.Item(74, 12).Value = 92742.8
$ws.Cells.Item(74, 13).Value = -101292.445
$ws.Cells.Item(74, 14).Value = -94490.8
$ws.Cells.Item(77, 8).Value = 97206.63
$ws.Cells.Item(77, 9).Value = 102166.445
$ws.Cells.Item(77, 10).Value = 92742.8
$ws.Cells.Item(77, 11).Value = 510832.225
$ws.Cells.Item(77, 12).Value = 463714
$ws.Cells.Item(77, 13).Value = -506464.225
$ws.Cells.Item(77, 14).Value = -472450
$ws.Cells.Item(97, 8).Value = 636795.7
$ws.Cells.Item(97, 9).Value = 927364.1
$ws.Cells.Item(97, 11).Value = 927364.1
$ws.Cells.Item(97, 13).Value = -926868.1
$ws.Cells.Item(102, 8).Value = 9263755
$ws.Cells.Item(102, 9).Value = 11908685
$ws.Cells.Item(102, 11).Value = 11908685
$ws.Cells.Item(102, 13).Value = -11907063
$ws.Cells.Item(122, 8).Value = 1161110.5
$ws.Cells.Item(122, 9).Value = 2787.1428
$ws.Cells.Item(122, 11).Value = 8361.428400000001
$ws.Cells.Item(122, 13).Value = -5911.428400000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3593.8096
$ws.Cells.Item(20, 9).Value = 2976.5
$ws.Cells.Item(20, 10).Value = 4416.8887
$ws.Cells.Item(20, 11).Value = 2976.5
$ws.Cells.Item(20, 12).Value = 4416.8887
$ws.Cells.Item(20, 13).Value = -2729.5
$ws.Cells.Item(20, 14).Value = -4910.8887
$ws.Cells.Item(31, 8).Value = 8266.666999999999
$ws.Cells.Item(31, 10).Value = 9900
$ws.Cells.Item(31, 12).Value = 9900
$ws.Cells.Item(31, 14).Value = -10404
$ws.Cells.Item(86, 8).Value = 5896506
$ws.Cells.Item(86, 9).Value = 7152528.5
$ws.Cells.Item(86, 10).Value = 35066.668
$ws.Cells.Item(86, 11).Value = 7152528.5
$ws.Cells.Item(86, 12).Value = 35066.668
$ws.Cells.Item(86, 13).Value = -7151405.5
$ws.Cells.Item(86, 14).Value = -37312.668
$ws.Cells.Item(89, 8).Value = 5896506
$ws.Cells.Item(89, 9).Value = 7152528.5
$ws.Cells.Item(89, 10).Value = 35066.668
$ws.Cells.Item(89, 11).Value = 35762642.5
$ws.Cells.Item(89, 12).Value = 175333.34
$ws.Cells.Item(89, 13).Value = -35757026.5
$ws.Cells.Item(89, 14).Value = -186565.34
$ws.Cells.Item(94, 8).Value = 2635330.8
$ws.Cells.Item(94, 9).Value = 3572827.5
$ws.Cells.Item(94, 10).Value = 10340.3
$ws.Cells.Item(94, 11).Value = 3572827.5
$ws.Cells.Item(94, 12).Value = 10340.3
$ws.Cells.Item(94, 13).Value = -3572376.5
$ws.Cells.Item(94, 14).Value = -11242.3
$ws.Cells.Item(105, 8).Value = 12500819
$ws.Cells.Item(105, 9).Value = 12500819
$ws.Cells.Item(105, 11).Value = 12500819
$ws.Cells.Item(105, 13).Value = -12499072

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 297.5
$ws.Cells.Item(22, 9).Value = 296.66666
$ws.Cells.Item(22, 11).Value = 296.66666
$ws.Cells.Item(22, 13).Value = 53.33334000000002
$ws.Cells.Item(31, 8).Value = 28633.572
$ws.Cells.Item(31, 9).Value = 1872.4
$ws.Cells.Item(31, 10).Value = 64315.133
$ws.Cells.Item(31, 11).Value = 1872.4
$ws.Cells.Item(31, 12).Value = 64315.133
$ws.Cells.Item(31, 13).Value = -1577.4
$ws.Cells.Item(31, 14).Value = -64905.133
$ws.Cells.Item(34, 8).Value = 28633.572
$ws.Cells.Item(34, 9).Value = 1872.4
$ws.Cells.Item(34, 10).Value = 64315.133
$ws.Cells.Item(34, 11).Value = 1872.4
$ws.Cells.Item(34, 12).Value = 64315.133
$ws.Cells.Item(34, 13).Value = -1670.4
$ws.Cells.Item(34, 14).Value = -64719.133

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 241.6923
$ws.Cells.Item(8, 9).Value = 241.6923
$ws.Cells.Item(8, 11).Value = 725.0769
$ws.Cells.Item(8, 13).Value = -586.0769
$ws.Cells.Item(56, 8).Value = 15630256
$ws.Cells.Item(56, 9).Value = 15630256
$ws.Cells.Item(56, 11).Value = 15630256
$ws.Cells.Item(56, 13).Value = -15629726
$ws.Cells.Item(125, 8).Value = 7573.5835
$ws.Cells.Item(125, 9).Value = 1000
$ws.Cells.Item(125, 11).Value = 3000
$ws.Cells.Item(125, 13).Value = 1920
$ws.Cells.Item(130, 8).Value = 2806.5454
$ws.Cells.Item(130, 9).Value = 2574.4
$ws.Cells.Item(130, 10).Value = 3000
$ws.Cells.Item(130, 11).Value = 7723.200000000001
$ws.Cells.Item(130, 12).Value = 9000
$ws.Cells.Item(130, 13).Value = -2703.200000000001
$ws.Cells.Item(130, 14).Value = -19040
$ws.Cells.Item(134, 8).Value = 3000
$ws.Cells.Item(134, 9).Value = 3000
$ws.Cells.Item(134, 11).Value = 9000
$ws.Cells.Item(134, 13).Value = -3930

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3581222.2
$ws.Cells.Item(126, 9).Value = 4548190
$ws.Cells.Item(126, 11).Value = 13644570
$ws.Cells.Item(126, 13).Value = -13642100

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1801.7307
$ws.Cells.Item(16, 9).Value = 1197.2222
$ws.Cells.Item(16, 10).Value = 3161.875
$ws.Cells.Item(16, 11).Value = 1197.2222
$ws.Cells.Item(16, 12).Value = 3161.875
$ws.Cells.Item(16, 13).Value = -1027.2222
$ws.Cells.Item(16, 14).Value = -3501.875
$ws.Cells.Item(22, 8).Value = 112274
$ws.Cells.Item(22, 9).Value = 222797.5
$ws.Cells.Item(22, 11).Value = 222797.5
$ws.Cells.Item(22, 13).Value = -222502.5
$ws.Cells.Item(27, 8).Value = 112274
$ws.Cells.Item(27, 9).Value = 222797.5
$ws.Cells.Item(27, 11).Value = 222797.5
$ws.Cells.Item(27, 13).Value = -222690.5
$ws.Cells.Item(55, 8).Value = 1502.3636
$ws.Cells.Item(55, 9).Value = 1722.2727
$ws.Cells.Item(55, 10).Value = 1282.4546
$ws.Cells.Item(55, 11).Value = 1722.2727
$ws.Cells.Item(55, 12).Value = 1282.4546
$ws.Cells.Item(55, 13).Value = -1549.2727
$ws.Cells.Item(55, 14).Value = -1628.4546
$ws.Cells.Item(132, 8).Value = 4778.8276
$ws.Cells.Item(132, 9).Value = 3745.9048
$ws.Cells.Item(132, 10).Value = 7490.25
$ws.Cells.Item(132, 11).Value = 11237.7144
$ws.Cells.Item(132, 12).Value = 22470.75
$ws.Cells.Item(132, 13).Value = -8707.714399999999
$ws.Cells.Item(132, 14).Value = -27530.75
$ws.Cells.Item(136, 8).Value = 59833.777
$ws.Cells.Item(136, 9).Value = 93596.27
$ws.Cells.Item(136, 10).Value = 6778.4287
$ws.Cells.Item(136, 11).Value = 280788.81
$ws.Cells.Item(136, 12).Value = 20335.2861
$ws.Cells.Item(136, 13).Value = -278238.81
$ws.Cells.Item(136, 14).Value = -25435.2861

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 5280
$ws.Cells.Item(96, 9).Value = 5136.2
$ws.Cells.Item(96, 11).Value = 5136.2
$ws.Cells.Item(96, 13).Value = -3763.2
$ws.Cells.Item(100, 8).Value = 1429.8182
$ws.Cells.Item(100, 9).Value = 1726.5
$ws.Cells.Item(100, 11).Value = 3453
$ws.Cells.Item(100, 13).Value = -2912
$ws.Cells.Item(126, 8).Value = 3632.75
$ws.Cells.Item(126, 9).Value = 3807
$ws.Cells.Item(126, 10).Value = 3388.8
$ws.Cells.Item(126, 11).Value = 11421
$ws.Cells.Item(126, 12).Value = 10166.4
$ws.Cells.Item(126, 13).Value = -8951
$ws.Cells.Item(126, 14).Value = -15106.4
$ws.Cells.Item(132, 8).Value = 18721908
$ws.Cells.Item(132, 9).Value = 23811398
$ws.Cells.Item(132, 11).Value = 71434194
$ws.Cells.Item(132, 13).Value = -71431664
$ws.Cells.Item(136, 8).Value = 2117.9736
$ws.Cells.Item(136, 9).Value = 1593.4375
$ws.Cells.Item(136, 11).Value = 4780.3125
$ws.Cells.Item(136, 13).Value = -2230.3125
